$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '27.031.50'
Set-TextValue 'E2' '  -0.28%  '
Set-TextValue 'D3' '1.831.22'
Set-TextValue 'D4' '1.006'
Set-TextValue 'E4' '  -0.53%  '
Set-TextValue 'D5' '312.41'
Set-TextValue 'E5' '  +0.14%  '
Set-TextValue 'E6' '  -0.51%  '
Set-TextValue 'D7' '0.4610'
Set-TextValue 'E7' '  -0.29%  '
Set-TextValue 'D8' '0.3708'
Set-TextValue 'E8' '  +1.82%  '
Set-TextValue 'D9' '0.07344'
Set-TextValue 'E9' '  +0.62%  '
Set-TextValue 'D10' '0.8757'
Set-TextValue 'E10' '  +0.69%  '
Set-TextValue 'D11' '0.07933'
Set-TextValue 'E11' '  +4.56%  '
Set-TextValue 'D12' '19.84'
Set-TextValue 'E12' '  -1.36%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.346'
Set-TextValue 'E13' '  +0.00%  '
Set-TextValue 'B14' 'Chainlink'
Set-TextValue 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '6.570'
Set-TextValue 'E14' '  +1.29%  '
Set-TextValue 'B15' 'WrappedEther'
Set-TextValue 'C15' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D15' '1.728.74'
Set-TextValue 'E15' '  -5.72%  '
Set-TextValue 'D16' '91.57'
Set-TextValue 'E16' '  -0.83%  '
Set-TextValue 'E17' '  -0.33%  '
Set-TextValue 'D18' '0.000008907'
Set-TextValue 'E18' '  +3.25%  '
Set-TextValue 'D19' '1.005'
Set-TextValue 'E19' '  -0.47%  '
Set-TextValue 'E20' '  +2.29%  '
Set-TextValue 'D21' '27.092.06'
Set-TextValue 'E21' '  -1.19%  '
Set-TextValue 'D22' '5.121'
Set-TextValue 'E22' '  -1.67%  '
Set-TextValue 'D23' '10.56'
Set-TextValue 'E23' '  +0.04%  '
Set-TextValue 'D24' '2.066.35'
Set-TextValue 'E24' '  -1.34%  '
Set-TextValue 'D25' '153.22'
Set-TextValue 'E25' '  +0.80%  '
Set-TextValue 'D26' '1.841'
Set-TextValue 'E26' '  -1.91%  '
Set-TextValue 'D27' '18.43'
Set-TextValue 'E27' '  +1.02%  '
Set-TextValue 'D28' '2.049'
Set-TextValue 'E28' '  -1.92%  '
Set-TextValue 'D29' '5.161'
Set-TextValue 'E29' '  +1.53%  '
Set-TextValue 'E30' '  -0.54%  '
Set-TextValue 'D31' '0.08909'
Set-TextValue 'E31' '  +0.03%  '
Set-TextValue 'D32' '2.964'
Set-TextValue 'E32' '  +0.09%  '
Set-TextValue 'D33' '0.7334'
Set-TextValue 'E33' '  +0.06%  '
Set-TextValue 'E34' '  -0.56%  '
Set-TextValue 'E35' '  -0.50%  '
Set-TextValue 'D36' '2.490'
Set-TextValue 'E36' '  +0.78%  '
Set-TextValue 'D37' '0.01952'
Set-TextValue 'E37' '  +1.92%  '
Set-TextValue 'E38' '  -0.20%  '
Set-TextValue 'D39' '0.05239'
Set-TextValue 'E39' '  -0.26%  '
Set-TextValue 'D40' '2.943'
Set-TextValue 'E40' '  +0.37%  '
Set-TextValue 'D41' '7.106'
Set-TextValue 'E41' '  -0.43%  '
Set-TextValue 'D42' '0.5165'
Set-TextValue 'E42' '  -0.61%  '
Set-TextValue 'B43' 'Frax'
Set-TextValue 'C43' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D43' '0.8948'
Set-TextValue 'E43' '  -11.50%  '
Set-TextValue 'B44' 'Algorand'
Set-TextValue 'C44' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.1628'
Set-TextValue 'E44' '  -0.13%  '
Set-TextValue 'B45' 'Aptos'
Set-TextValue 'C45' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D45' '8.220'
Set-TextValue 'E45' '  -0.50%  '
Set-TextValue 'B46' 'Decentraland'
Set-TextValue 'C46' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D46' '0.4842'
Set-TextValue 'E46' '  -0.54%  '
Set-TextValue 'D47' '10.23'
Set-TextValue 'E47' '  +0.44%  '
Set-TextValue 'B48' 'PaxDollar'
Set-TextValue 'C48' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D48' '1.005'
Set-TextValue 'E48' '  -0.53%  '
Set-TextValue 'B49' 'Quant'
Set-TextValue 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D49' '102.31'
Set-TextValue 'E49' '  -1.26%  '
Set-TextValue 'B50' 'NEARProtocol'
Set-TextValue 'C50' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D50' '1.632'
Set-TextValue 'E50' '  +0.00%  '
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.06201'
Set-TextValue 'E51' '  -0.98%  '
